# Apply league base update (17-02-2024 11:11) to "Denmark Division 2" sheet.
# For several fixtures the data rows (columns B..AC) were rearranged
# (row identity/order in column A is kept; the match data that used to sit
# in one row now sits in another row of the same small cluster).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29)

# --- rows [24, 25] ---
$snap24 = @()
foreach ($c in $cols) { $snap24 += $ws.Cells.Item(24, $c).Value2 }
$snap25 = @()
foreach ($c in $cols) { $snap25 += $ws.Cells.Item(25, $c).Value2 }
for ($i = 0; $i -lt $cols.Count; $i++) { $ws.Cells.Item(24, $cols[$i]).Value2 = $snap25[$i] }
for ($i = 0; $i -lt $cols.Count; $i++) { $ws.Cells.Item(25, $cols[$i]).Value2 = $snap24[$i] }

# --- rows [26, 27] ---
$snap26 = @()
foreach ($c in $cols) { $snap26 += $ws.Cells.Item(26, $c).Value2 }
$snap27 = @()
foreach ($c in $cols) { $snap27 += $ws.Cells.Item(27, $c).Value2 }
for ($i = 0; $i -lt $cols.Count; $i++) { $ws.Cells.Item(26, $cols[$i]).Value2 = $snap27[$i] }
for ($i = 0; $i -lt $cols.Count; $i++) { $ws.Cells.Item(27, $cols[$i]).Value2 = $snap26[$i] }

# --- rows [31, 32] ---
$snap31 = @()
foreach ($c in $cols) { $snap31 += $ws.Cells.Item(31, $c).Value2 }
$snap32 = @()
foreach ($c in $cols) { $snap32 += $ws.Cells.Item(32, $c).Value2 }
for ($i = 0; $i -lt $cols.Count; $i++) { $ws.Cells.Item(31, $cols[$i]).Value2 = $snap32[$i] }
for ($i = 0; $i -lt $cols.Count; $i++) { $ws.Cells.Item(32, $cols[$i]).Value2 = $snap31[$i] }

# --- rows [65, 66] ---
$snap65 = @()
foreach ($c in $cols) { $snap65 += $ws.Cells.Item(65, $c).Value2 }
$snap66 = @()
foreach ($c in $cols) { $snap66 += $ws.Cells.Item(66, $c).Value2 }
for ($i = 0; $i -lt $cols.Count; $i++) { $ws.Cells.Item(65, $cols[$i]).Value2 = $snap66[$i] }
for ($i = 0; $i -lt $cols.Count; $i++) { $ws.Cells.Item(66, $cols[$i]).Value2 = $snap65[$i] }

# --- rows [90, 91] ---
$snap90 = @()
foreach ($c in $cols) { $snap90 += $ws.Cells.Item(90, $c).Value2 }
$snap91 = @()
foreach ($c in $cols) { $snap91 += $ws.Cells.Item(91, $c).Value2 }
for ($i = 0; $i -lt $cols.Count; $i++) { $ws.Cells.Item(90, $cols[$i]).Value2 = $snap91[$i] }
for ($i = 0; $i -lt $cols.Count; $i++) { $ws.Cells.Item(91, $cols[$i]).Value2 = $snap90[$i] }

# --- rows [95, 96, 97] ---
$snap95 = @()
foreach ($c in $cols) { $snap95 += $ws.Cells.Item(95, $c).Value2 }
$snap96 = @()
foreach ($c in $cols) { $snap96 += $ws.Cells.Item(96, $c).Value2 }
$snap97 = @()
foreach ($c in $cols) { $snap97 += $ws.Cells.Item(97, $c).Value2 }
for ($i = 0; $i -lt $cols.Count; $i++) { $ws.Cells.Item(95, $cols[$i]).Value2 = $snap97[$i] }
for ($i = 0; $i -lt $cols.Count; $i++) { $ws.Cells.Item(96, $cols[$i]).Value2 = $snap95[$i] }
for ($i = 0; $i -lt $cols.Count; $i++) { $ws.Cells.Item(97, $cols[$i]).Value2 = $snap96[$i] }

# --- rows [149, 150, 151] ---
$snap149 = @()
foreach ($c in $cols) { $snap149 += $ws.Cells.Item(149, $c).Value2 }
$snap150 = @()
foreach ($c in $cols) { $snap150 += $ws.Cells.Item(150, $c).Value2 }
$snap151 = @()
foreach ($c in $cols) { $snap151 += $ws.Cells.Item(151, $c).Value2 }
for ($i = 0; $i -lt $cols.Count; $i++) { $ws.Cells.Item(149, $cols[$i]).Value2 = $snap151[$i] }
for ($i = 0; $i -lt $cols.Count; $i++) { $ws.Cells.Item(150, $cols[$i]).Value2 = $snap149[$i] }
for ($i = 0; $i -lt $cols.Count; $i++) { $ws.Cells.Item(151, $cols[$i]).Value2 = $snap150[$i] }
